$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.448.86"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "'1.725.67"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'244.31"
$ws.Range("E5").Value = "  +2.11%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4802"
$ws.Range("E7").Value = "  +2.73%  "
$ws.Range("D8").Value = "'0.2683"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("D9").Value = "'0.06221"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").Value = "'1.730.62"
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("D11").Value = "'0.07120"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").Value = "'15.69"
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("D13").Value = "'0.6159"
$ws.Range("E13").Value = "  +5.30%  "
$ws.Range("D14").Value = "'4.543"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "'77.11"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'26.453.81"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "'11.71"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "'1.953.08"
$ws.Range("E21").Value = "  +2.95%  "
$ws.Range("D22").Value = "'4.537"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'8.895"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Value = "'5.308"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'136.26"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("D26").Value = "'15.38"
$ws.Range("E26").Value = "  +2.02%  "
$ws.Range("D27").Value = "'1.794"
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("D28").Value = "'1.411"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "'106.90"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").Value = "'3.970"
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("D31").Value = "'0.08038"
$ws.Range("E31").Value = "  +3.59%  "
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'0.04548"
$ws.Range("E33").Value = "  +3.95%  "
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").Value = "'0.6362"
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").Value = "'0.9878"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").Value = "'0.9356"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").Value = "'1.992"
$ws.Range("E38").Value = "  +5.42%  "
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").Value = "'107.44"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.401"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "'0.01496"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("E43").Value = "  +10.99%  "
$ws.Range("D44").Value = "'0.3908"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("D45").Value = "'6.982"
$ws.Range("E45").Value = "  +12.37%  "
$ws.Range("D46").Value = "'0.1190"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("D47").Value = "'0.05319"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'30.90"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.850"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").Value = "'1.266"
$ws.Range("E50").Value = "  +4.24%  "

# Reset style on forced-text numeric cells to avoid leaving quotePrefix formatting
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
